$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9437097907066345
$ws.Range("B1").Value = 1.998593091964722
$ws.Range("C1").Value = 7.728731155395508
$ws.Range("D1").Value = 2.69521951675415
$ws.Range("E1").Value = 1.092660784721375
